# "lock in current version" - refresh the word/image/category cue lists
# (columns A, B, C for rows 2:49) with the new curated set of German verbs
# and dog/flower image cues.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "wenden",
    "drücken",
    "nerven",
    "rufen",
    "schleppen",
    "ehren",
    "enden",
    "lügen",
    "kommen",
    "opfern",
    "treten",
    "machen",
    "schalten",
    "tropfen",
    "schützen",
    "klagen",
    "quälen",
    "kümmern",
    "sparen",
    "packen",
    "fragen",
    "hören",
    "sehen",
    "triefen",
    "drohen",
    "bilden",
    "zögern",
    "orten",
    "achten",
    "dringen",
    "weigern",
    "knien",
    "spielen",
    "ächzen",
    "hoffen",
    "süßen",
    "kosten",
    "wachsen",
    "fügen",
    "dauern",
    "fordern",
    "ändern",
    "stören",
    "rasen",
    "tauschen",
    "bremsen",
    "heben",
    "ärgern"
)

$images = @(
    "none",
    "flower/flower004.jpg",
    "dog/dog013.jpg",
    "none",
    "dog/dog025.jpg",
    "dog/dog027.jpg",
    "none",
    "dog/dog002.jpg",
    "dog/dog019.jpg",
    "none",
    "dog/dog016.jpg",
    "dog/dog030.jpg",
    "none",
    "flower/flower011.jpg",
    "flower/flower023.jpg",
    "none",
    "dog/dog029.jpg",
    "flower/flower032.jpg",
    "none",
    "flower/flower012.jpg",
    "flower/flower031.jpg",
    "none",
    "dog/dog020.jpg",
    "dog/dog012.jpg",
    "none",
    "flower/flower020.jpg",
    "dog/dog009.jpg",
    "none",
    "dog/dog021.jpg",
    "dog/dog015.jpg",
    "none",
    "flower/flower022.jpg",
    "flower/flower008.jpg",
    "none",
    "flower/flower010.jpg",
    "flower/flower001.jpg",
    "none",
    "flower/flower029.jpg",
    "flower/flower019.jpg",
    "none",
    "dog/dog022.jpg",
    "dog/dog031.jpg",
    "none",
    "dog/dog001.jpg",
    "flower/flower030.jpg",
    "none",
    "flower/flower007.jpg",
    "flower/flower003.jpg"
)

$categories = @(
    "none",
    "flower",
    "dog",
    "none",
    "dog",
    "dog",
    "none",
    "dog",
    "dog",
    "none",
    "dog",
    "dog",
    "none",
    "flower",
    "flower",
    "none",
    "dog",
    "flower",
    "none",
    "flower",
    "flower",
    "none",
    "dog",
    "dog",
    "none",
    "flower",
    "dog",
    "none",
    "dog",
    "dog",
    "none",
    "flower",
    "flower",
    "none",
    "flower",
    "flower",
    "none",
    "flower",
    "flower",
    "none",
    "dog",
    "dog",
    "none",
    "dog",
    "flower",
    "none",
    "flower",
    "flower"
)

for ($i = 0; $i -lt $words.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $words[$i]
}

for ($i = 0; $i -lt $images.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $images[$i]
}

for ($i = 0; $i -lt $categories.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $categories[$i]
}
